$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MAX (Q2) to a formula literal "=1591" (cached value 1591)
$ws.Range("Q2").Formula = "=1591"

# Update MIN (Q5) to -18
$ws.Range("Q5").Value = -18

# Move the selection/active cell to O20 (no data there, just cursor position)
$ws.Range("O20").Select()
